$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-11 Sunday" "2024-08-12 Monday"

Replace-Text "204÷4=51, 0" "861÷2=430, 1"
Replace-Text "589÷2=294, 1" "339÷4=84, 3"
Replace-Text "834÷9=92, 6" "889÷6=148, 1"
Replace-Text "115÷4=28, 3" "820÷7=117, 1"
Replace-Text "672÷7=96, 0" "476÷4=119, 0"

Replace-Text "909÷8=113, 5" "267÷2=133, 1"
Replace-Text "508÷4=127, 0" "977÷4=244, 1"
Replace-Text "662÷3=220, 2" "407÷3=135, 2"
Replace-Text "825÷9=91, 6" "423÷7=60, 3"
Replace-Text "155÷7=22, 1" "165÷9=18, 3"

Replace-Text "569÷8=71, 1" "483÷2=241, 1"
Replace-Text "896÷7=128, 0" "585÷8=73, 1"
Replace-Text "440÷2=220, 0" "290÷7=41, 3"
Replace-Text "213÷9=23, 6" "205÷6=34, 1"
Replace-Text "269÷7=38, 3" "472÷6=78, 4"

Replace-Text "134÷9=14, 8" "900÷8=112, 4"
Replace-Text "666÷9=74, 0" "803÷7=114, 5"
Replace-Text "785÷2=392, 1" "783÷6=130, 3"
Replace-Text "325÷2=162, 1" "370÷3=123, 1"
Replace-Text "866÷8=108, 2" "943÷4=235, 3"

Replace-Text "882÷6=147, 0" "450÷9=50, 0"
Replace-Text "252÷7=36, 0" "246÷2=123, 0"
Replace-Text "846÷9=94, 0" "800÷4=200, 0"
Replace-Text "554÷3=184, 2" "579÷7=82, 5"
Replace-Text "425÷7=60, 5" "890÷3=296, 2"
